$d = $word.ActiveDocument

# --- 1. Insert the new block of paragraphs/table at the end of the document ---
# (Done before touching $d.Tables, since indexing into Tables disturbs later
#  Paragraphs.Item(...) lookups in this runtime.)
# The new content is inserted right before the very last (empty) paragraph of
# the document, i.e. between the two trailing empty paragraphs.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertRange = $lastPara.Range
$insertRange.Collapse(1)

$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p/><w:p/><w:p/><w:p/><w:p/><w:p/><w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0"/></w:tblPr><w:tblGrid><w:gridCol w:w="2718"/><w:gridCol w:w="791"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="2718" w:type="dxa"/><w:shd w:val="clear" w:color="auto" w:fill="D9D9D9" w:themeFill="background1" w:themeFillShade="D9"/></w:tcPr><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Materi</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="791" w:type="dxa"/><w:shd w:val="clear" w:color="auto" w:fill="D9D9D9" w:themeFill="background1" w:themeFillShade="D9"/></w:tcPr><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Nilai </w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="2718" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Mengerjakan Video 1-19</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="791" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/></w:rPr><w:t>80</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Saya Sudah Belajar dan </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>M</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>engerti</w:t></w:r></w:p><w:p><w:r><w:t>1. Saya Sudah Mengerti Cara Membuat Folder CI-4</w:t></w:r></w:p><w:p><w:r><w:t>2.Saya Belajar Kembali CI-4, dan Berlatih</w:t></w:r><w:r><w:tab/></w:r></w:p><w:p><w:r><w:t>3</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Saya Belum Mengerti</w:t></w:r></w:p><w:p><w:r><w:t>1.Saya Belum Mengerti Bagaimana Cara Memahami CI-4 Dengan Cepat</w:t></w:r></w:p><w:p><w:r><w:t>2.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">3. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertRange.InsertXML($xml)

# --- 2. Resize the second (pre-existing) table's columns (3969/1083 dxa -> 2718/791 dxa) ---
# 1 point = 20 dxa, so 2718 dxa = 135.9 pt, 791 dxa = 39.55 pt
$t2 = $d.Tables.Item(2)
$t2.Columns.Item(1).Width = 135.9
$t2.Columns.Item(2).Width = 39.55
